# Apply the "Updated database schema diagram" edit to slide 2 of the
# DataModel deck.
#
#  - Table "RuleHasAction" (shape 9): the data type of the "action"
#    column changes from "int" to "text".
#  - Table "RuleHasAssignor" (shape 17): the "party_id"/"int" and
#    "rule_id"/"int" rows become bold + italic (matching the styling
#    already used on the equivalent rows of the other association
#    tables on this slide).
#  - Table "RuleHasAssignee" (shape 18): same bold/italic treatment as
#    RuleHasAssignor.

$p  = $ppt.ActivePresentation
$s2 = $p.Slides.Item(2)

# --- RuleHasAction: action's type int -> text -----------------------
$ruleHasAction = $s2.Shapes.Item(9)
$actionTypeCell = $ruleHasAction.Table.Cell(3, 2)
$actionTypeCell.Shape.TextFrame.TextRange.Text = "text"

# --- RuleHasAssignor: bold + italic the FK rows ----------------------
$ruleHasAssignor = $s2.Shapes.Item(17)
for ($r = 2; $r -le 3; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $tr = $ruleHasAssignor.Table.Cell($r, $c).Shape.TextFrame.TextRange
        $tr.Font.Bold = -1
        $tr.Font.Italic = -1
    }
}

# --- RuleHasAssignee: bold + italic the FK rows ----------------------
$ruleHasAssignee = $s2.Shapes.Item(18)
for ($r = 2; $r -le 3; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $tr = $ruleHasAssignee.Table.Cell($r, $c).Shape.TextFrame.TextRange
        $tr.Font.Bold = -1
        $tr.Font.Italic = -1
    }
}
